# Add the "OutputPaths" sheet (output path definitions) after the existing
# "Scenarios" sheet, add a new "OutputPathsIds" column to the Scenarios
# sheet, and link the second scenario row to the new output paths.

$wb = $excel.ActiveWorkbook
$scenarios = $wb.Worksheets.Item(1)

# New sheet goes after the last existing sheet ("Scenarios").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$outputPaths = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$outputPaths.Name = "OutputPaths"

# New header cell on Scenarios (column L) - written first so the new shared
# string "OutputPathsIds" gets allocated before the OutputPaths sheet values.
$scenarios.Range("L1").Value = "OutputPathsIds"

# OutputPaths sheet header + data.
$outputPaths.Range("A1").Value = "OutputPathId"
$outputPaths.Range("B1").Value = "OutputPath"
$outputPaths.Range("A2").Value = "Aciclovir_PVB"
$outputPaths.Range("A3").Value = "Aciclovir_fat_cell"
$outputPaths.Range("B2").Value = "Organism|PeripheralVenousBlood|Aciclovir|Plasma (Peripheral Venous Blood)"
$outputPaths.Range("B3").Value = "Organism|Fat|Intracellular|Aciclovir|Concentration in container"

# Header row is bold, like the rest of the workbook's header rows.
$outputPaths.Range("A1:B1").Font.Bold = $true

# Scenarios row 3 (TestScenario2) references both new output paths.
$scenarios.Range("L3").Value = "Aciclovir_PVB, Aciclovir_fat_cell"

# Column widths (best-fit-like sizing for the new columns).
$scenarios.Columns.Item(12).ColumnWidth = 14.86
$outputPaths.Columns.Item(1).ColumnWidth = 16.29
$outputPaths.Columns.Item(2).ColumnWidth = 11.29

# Selection state: OutputPaths ends up selected at A3, Scenarios (the
# active sheet) ends up selected at L3.
$outputPaths.Range("A3").Select() | Out-Null
$scenarios.Range("L3").Select() | Out-Null
